$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: fix property_category from "land" to "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (car) sheet: fix property_category from "land" to "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
